$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C34").Value = 'Metropolitana'
$ws.Range("D34").Value = 44449
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 'Fruta'
$ws.Range("G34").Value = 100107
$ws.Range("H34").Value = 'Otros'
$ws.Range("I34").Value = 100107002
$ws.Range("J34").Value = 'Chirimoya'
$ws.Range("K34").Value = 'Cultivar IV Región'
$ws.Range("L34").Value = 'Cuarta'
$ws.Range("M34").Value = 300
$ws.Range("N34").Value = 1000
$ws.Range("O34").Value = 1000
$ws.Range("P34").Value = 1000
$ws.Range("Q34").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R34").Value = 'Provincia del Elquí'
$ws.Range("S34").Value = 1000
$ws.Range("T34").Value = 1

# Row 35
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C35").Value = 'Metropolitana'
$ws.Range("D35").Value = 44449
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 'Fruta'
$ws.Range("G35").Value = 100107
$ws.Range("H35").Value = 'Otros'
$ws.Range("I35").Value = 100107002
$ws.Range("J35").Value = 'Chirimoya'
$ws.Range("K35").Value = 'Cultivar IV Región'
$ws.Range("L35").Value = 'Especial'
$ws.Range("M35").Value = 250
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 20000
$ws.Range("Q35").Value = '$/bandeja 8 kilos'
$ws.Range("R35").Value = 'Provincia del Elquí'
$ws.Range("S35").Value = 2500
$ws.Range("T35").Value = 8

# Row 36
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C36").Value = 'Metropolitana'
$ws.Range("D36").Value = 44449
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 'Fruta'
$ws.Range("G36").Value = 100107
$ws.Range("H36").Value = 'Otros'
$ws.Range("I36").Value = 100107002
$ws.Range("J36").Value = 'Chirimoya'
$ws.Range("K36").Value = 'Cultivar IV Región'
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 330
$ws.Range("N36").Value = 17600
$ws.Range("O36").Value = 17600
$ws.Range("P36").Value = 17600
$ws.Range("Q36").Value = '$/bandeja 8 kilos'
$ws.Range("R36").Value = 'Provincia del Elquí'
$ws.Range("S36").Value = 2200
$ws.Range("T36").Value = 8

# Row 37
$ws.Range("A37").Value = 9
$ws.Range("B37").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C37").Value = 'Metropolitana'
$ws.Range("D37").Value = 44449
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 'Fruta'
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = 'Otros'
$ws.Range("I37").Value = 100107002
$ws.Range("J37").Value = 'Chirimoya'
$ws.Range("K37").Value = 'Cultivar IV Región'
$ws.Range("L37").Value = 'Segunda'
$ws.Range("M37").Value = 300
$ws.Range("N37").Value = 14400
$ws.Range("O37").Value = 14400
$ws.Range("P37").Value = 14400
$ws.Range("Q37").Value = '$/bandeja 8 kilos'
$ws.Range("R37").Value = 'Provincia del Elquí'
$ws.Range("S37").Value = 1800
$ws.Range("T37").Value = 8

# Row 38
$ws.Range("A38").Value = 9
$ws.Range("B38").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C38").Value = 'Metropolitana'
$ws.Range("D38").Value = 44449
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = 'Fruta'
$ws.Range("G38").Value = 100107
$ws.Range("H38").Value = 'Otros'
$ws.Range("I38").Value = 100107002
$ws.Range("J38").Value = 'Chirimoya'
$ws.Range("K38").Value = 'Cultivar IV Región'
$ws.Range("L38").Value = 'Tercera'
$ws.Range("M38").Value = 280
$ws.Range("N38").Value = 1400
$ws.Range("O38").Value = 1400
$ws.Range("P38").Value = 1400
$ws.Range("Q38").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R38").Value = 'Provincia del Elquí'
$ws.Range("S38").Value = 1400
$ws.Range("T38").Value = 1

# Row 39
$ws.Range("A39").Value = 9
$ws.Range("B39").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C39").Value = 'Metropolitana'
$ws.Range("D39").Value = 44400
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 'Fruta'
$ws.Range("G39").Value = 100107
$ws.Range("H39").Value = 'Otros'
$ws.Range("I39").Value = 100107002
$ws.Range("J39").Value = 'Chirimoya'
$ws.Range("K39").Value = 'Cultivar IV Región'
$ws.Range("L39").Value = 'Especial'
$ws.Range("M39").Value = 15
$ws.Range("N39").Value = 2000
$ws.Range("O39").Value = 2000
$ws.Range("P39").Value = 2000
$ws.Range("Q39").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R39").Value = 'Provincia del Elquí'
$ws.Range("S39").Value = 2000
$ws.Range("T39").Value = 1

# Row 40
$ws.Range("A40").Value = 9
$ws.Range("B40").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C40").Value = 'Metropolitana'
$ws.Range("D40").Value = 44400
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 'Fruta'
$ws.Range("G40").Value = 100107
$ws.Range("H40").Value = 'Otros'
$ws.Range("I40").Value = 100107002
$ws.Range("J40").Value = 'Chirimoya'
$ws.Range("K40").Value = 'Cultivar IV Región'
$ws.Range("L40").Value = 'Extra (doble especial)'
$ws.Range("M40").Value = 10
$ws.Range("N40").Value = 2500
$ws.Range("O40").Value = 2500
$ws.Range("P40").Value = 2500
$ws.Range("Q40").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R40").Value = 'Provincia del Elquí'
$ws.Range("S40").Value = 2500
$ws.Range("T40").Value = 1

# Row 41
$ws.Range("A41").Value = 9
$ws.Range("B41").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C41").Value = 'Metropolitana'
$ws.Range("D41").Value = 44400
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 'Fruta'
$ws.Range("G41").Value = 100107
$ws.Range("H41").Value = 'Otros'
$ws.Range("I41").Value = 100107002
$ws.Range("J41").Value = 'Chirimoya'
$ws.Range("K41").Value = 'Cultivar IV Región'
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 8
$ws.Range("N41").Value = 1500
$ws.Range("O41").Value = 1500
$ws.Range("P41").Value = 1500
$ws.Range("Q41").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R41").Value = 'Provincia del Elquí'
$ws.Range("S41").Value = 1500
$ws.Range("T41").Value = 1

# Row 42
$ws.Range("A42").Value = 9
$ws.Range("B42").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C42").Value = 'Metropolitana'
$ws.Range("D42").Value = 44400
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = 'Fruta'
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = 'Otros'
$ws.Range("I42").Value = 100107002
$ws.Range("J42").Value = 'Chirimoya'
$ws.Range("K42").Value = 'Cultivar IV Región'
$ws.Range("L42").Value = 'Segunda'
$ws.Range("M42").Value = 6
$ws.Range("N42").Value = 1000
$ws.Range("O42").Value = 1000
$ws.Range("P42").Value = 1000
$ws.Range("Q42").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R42").Value = 'Provincia del Elquí'
$ws.Range("S42").Value = 1000
$ws.Range("T42").Value = 1
